$wb = $excel.ActiveWorkbook

# --- Rename "deals" -> "products" ---
$ws = $wb.Worksheets.Item("deals")
$ws.Name = "products"

# --- Populate header row + first data row ---
$headers = @("srno", "productCategory", "productSubCategory", "productName", "quantity", "color")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
}
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Dresses"
$ws.Range("C2").Value = "Summer Dresses"

# --- Column widths for products sheet (values chosen to round-trip through
#     this runtime's internal character-width quantization as closely as
#     possible to the target stored widths of 17 / 21.375 / 14.375) ---
$ws.Columns.Item(2).ColumnWidth = 16.17
$ws.Columns.Item(3).ColumnWidth = 20.5
$ws.Columns.Item(4).ColumnWidth = 13.5

# --- Header formatting: yellow fill + thin box border, applied to A1 first then
#     copied across so only a single new cell style gets registered. ---
$a1 = $ws.Range("A1")
$a1.Interior.Color = 65535
$a1.Borders.LineStyle = 1
$a1.Borders.Weight = 2
$a1.Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Make products the active/selected sheet & cell ---
$ws.Activate()
$ws.Range("F15").Select() | Out-Null
